# Scheduled runner update: refresh cached Universalis market-board price
# snapshots (currentAveragePrice / NQ / HQ) and the resulting Leve profit
# columns (H:N) across the Ultima_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1652.8889
$ws.Range("I40").Value = 1508
$ws.Range("J40").Value = 1834
$ws.Range("K40").Value = 1508
$ws.Range("L40").Value = 1834
$ws.Range("M40").Value = -1333
$ws.Range("N40").Value = -2184
$ws.Range("H64").Value = 12348568
$ws.Range("I64").Value = 27780178
$ws.Range("J64").Value = 3280
$ws.Range("K64").Value = 27780178
$ws.Range("L64").Value = 3280
$ws.Range("M64").Value = -27779930
$ws.Range("N64").Value = -3776
$ws.Range("H67").Value = 12348568
$ws.Range("I67").Value = 27780178
$ws.Range("J67").Value = 3280
$ws.Range("K67").Value = 27780178
$ws.Range("L67").Value = 3280
$ws.Range("M67").Value = -27779320
$ws.Range("N67").Value = -4996
$ws.Range("H116").Value = 2842.25
$ws.Range("I116").Value = 2650.875
$ws.Range("K116").Value = 2650.875
$ws.Range("M116").Value = 791.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 36367.332
$ws.Range("J48").Value = 36367.332
$ws.Range("L48").Value = 36367.332
$ws.Range("N48").Value = -37319.332
$ws.Range("H62").Value = 3440.5557
$ws.Range("I62").Value = 3573.75
$ws.Range("J62").Value = 2375
$ws.Range("K62").Value = 3573.75
$ws.Range("L62").Value = 2375
$ws.Range("M62").Value = -2949.75
$ws.Range("N62").Value = -3623
$ws.Range("H65").Value = 3440.5557
$ws.Range("I65").Value = 3573.75
$ws.Range("J65").Value = 2375
$ws.Range("K65").Value = 17868.75
$ws.Range("L65").Value = 11875
$ws.Range("M65").Value = -14748.75
$ws.Range("N65").Value = -18115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1682.4736
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 2578.818
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 7736.454000000001
$ws.Range("M22").Value = -1181
$ws.Range("N22").Value = -8074.454000000001
$ws.Range("H27").Value = 1682.4736
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 2578.818
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 7736.454000000001
$ws.Range("M27").Value = -1248
$ws.Range("N27").Value = -7940.454000000001
$ws.Range("H63").Value = 4729.7
$ws.Range("I63").Value = 2659.4
$ws.Range("J63").Value = 6800
$ws.Range("K63").Value = 7978.200000000001
$ws.Range("L63").Value = 20400
$ws.Range("M63").Value = -7229.200000000001
$ws.Range("N63").Value = -21898
$ws.Range("H66").Value = 4729.7
$ws.Range("I66").Value = 2659.4
$ws.Range("J66").Value = 6800
$ws.Range("K66").Value = 23934.6
$ws.Range("L66").Value = 61200
$ws.Range("M66").Value = -20190.6
$ws.Range("N66").Value = -68688
$ws.Range("H76").Value = 4000
$ws.Range("J76").Value = 4000
$ws.Range("L76").Value = 12000
$ws.Range("N76").Value = -12766
$ws.Range("H79").Value = 4000
$ws.Range("J79").Value = 4000
$ws.Range("L79").Value = 12000
$ws.Range("N79").Value = -14652
$ws.Range("H81").Value = 1137.2858
$ws.Range("I81").Value = 547.875
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 1643.625
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -520.625
$ws.Range("N81").Value = -6746
$ws.Range("H84").Value = 1137.2858
$ws.Range("I84").Value = 547.875
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 4930.875
$ws.Range("L84").Value = 13500
$ws.Range("M84").Value = 685.125
$ws.Range("N84").Value = -24732
$ws.Range("H86").Value = 1512.381
$ws.Range("I86").Value = 958.8889
$ws.Range("J86").Value = 1927.5
$ws.Range("K86").Value = 2876.6667
$ws.Range("L86").Value = 5782.5
$ws.Range("M86").Value = -1690.6667
$ws.Range("N86").Value = -8154.5
$ws.Range("H87").Value = 15490
$ws.Range("I87").Value = 8000
$ws.Range("K87").Value = 24000
$ws.Range("M87").Value = -22752
$ws.Range("H89").Value = 1512.381
$ws.Range("I89").Value = 958.8889
$ws.Range("J89").Value = 1927.5
$ws.Range("K89").Value = 8630.000100000001
$ws.Range("L89").Value = 17347.5
$ws.Range("M89").Value = -2702.000100000001
$ws.Range("N89").Value = -29203.5
$ws.Range("H90").Value = 15490
$ws.Range("I90").Value = 8000
$ws.Range("K90").Value = 72000
$ws.Range("M90").Value = -65760
$ws.Range("H92").Value = 1455.8889
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 1512.875
$ws.Range("K92").Value = 3000
$ws.Range("L92").Value = 4538.625
$ws.Range("M92").Value = -1752
$ws.Range("N92").Value = -7034.625
$ws.Range("H94").Value = 6250
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 12000
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 36000
$ws.Range("M94").Value = -824
$ws.Range("N94").Value = -37352
$ws.Range("H107").Value = 924.38464
$ws.Range("I107").Value = 165
$ws.Range("J107").Value = 987.6667
$ws.Range("K107").Value = 495
$ws.Range("L107").Value = 2963.0001
$ws.Range("M107").Value = 1425
$ws.Range("N107").Value = -6803.0001
$ws.Range("H109").Value = 2649.05
$ws.Range("I109").Value = 1256.2
$ws.Range("J109").Value = 3113.3333
$ws.Range("K109").Value = 3768.6
$ws.Range("L109").Value = 9339.999899999999
$ws.Range("M109").Value = -2728.6
$ws.Range("N109").Value = -11419.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3901.8572
$ws.Range("I126").Value = 2453.3845
$ws.Range("J126").Value = 5157.2
$ws.Range("K126").Value = 7360.1535
$ws.Range("L126").Value = 15471.6
$ws.Range("M126").Value = -4890.1535
$ws.Range("N126").Value = -20411.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4843.853
$ws.Range("I7").Value = 5724.25
$ws.Range("J7").Value = 4363.636
$ws.Range("K7").Value = 5724.25
$ws.Range("L7").Value = 4363.636
$ws.Range("M7").Value = -5612.25
$ws.Range("N7").Value = -4587.636
$ws.Range("H122").Value = 5817.433
$ws.Range("I122").Value = 5576.278
$ws.Range("J122").Value = 6179.1665
$ws.Range("K122").Value = 16728.834
$ws.Range("L122").Value = 18537.4995
$ws.Range("M122").Value = -14278.834
$ws.Range("N122").Value = -23437.4995
$ws.Range("H126").Value = 4843.853
$ws.Range("I126").Value = 5724.25
$ws.Range("J126").Value = 4363.636
$ws.Range("K126").Value = 17172.75
$ws.Range("L126").Value = 13090.908
$ws.Range("M126").Value = -14702.75
$ws.Range("N126").Value = -18030.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 12000
$ws.Range("J93").Value = 12000
$ws.Range("L93").Value = 12000
$ws.Range("N93").Value = -16992
$ws.Range("H113").Value = 1682.2354
$ws.Range("I113").Value = 325
$ws.Range("J113").Value = 2888.6667
$ws.Range("K113").Value = 975
$ws.Range("L113").Value = 8666.000100000001
$ws.Range("M113").Value = 1195
$ws.Range("N113").Value = -13006.0001
